$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.588.07'
$ws.Range('E2').Value = '  -1.90%  '
$ws.Range('D3').Value = '2.889.34'
$ws.Range('E3').Value = '  -1.77%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = "'565.50"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.50%  '
$ws.Range('D6').Value = "'142.49"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.17%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'0.499"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.82%  '
$ws.Range('D9').Value = '2.885.16'
$ws.Range('E9').Value = '  -1.79%  '
$ws.Range('E10').Value = '  -2.00%  '
$ws.Range('D11').Value = "'0.145"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.98%  '
$ws.Range('E12').Value = '  -2.45%  '
$ws.Range('D13').Value = "'0.0000230"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.73%  '
$ws.Range('D14').Value = "'31.55"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.74%  '
$ws.Range('D15').Value = "'0.126"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').Value = '3.371.70'
$ws.Range('D17').Value = '61.597.24'
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('D18').Value = '2.900.04'
$ws.Range('E18').Value = '  -1.50%  '
$ws.Range('D19').Value = "'6.48"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.73%  '
$ws.Range('D20').Value = "'428.08"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.25%  '
$ws.Range('D21').Value = "'12.96"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.25%  '
$ws.Range('D22').Value = "'0.649"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.09%  '
$ws.Range('D23').Value = "'6.78"
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Value = "'78.76"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.41%  '
$ws.Range('D25').Value = "'11.83"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.26%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = "'9.89"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -11.31%  '
$ws.Range('E28').Value = '  -5.82%  '
$ws.Range('D29').Value = "'0.0000105"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +7.58%  '
$ws.Range('D30').Value = "'6.99"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.10%  '
$ws.Range('D31').Value = "'2.48"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.56%  '
$ws.Range('D32').Value = "'2.02"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -9.38%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').Value = "'0.105"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.11%  '
$ws.Range('D35').Value = "'25.38"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.36%  '
$ws.Range('D36').Value = "'0.958"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.32%  '
$ws.Range('D37').Value = "'5.33"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.86%  '
$ws.Range('D38').Value = "'48.77"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.72%  '
$ws.Range('E39').Value = '  -4.76%  '
$ws.Range('D40').Value = "'2.78"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.84%  '
$ws.Range('D41').Value = "'8.13"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.54%  '
$ws.Range('E42').Value = '  -4.28%  '
$ws.Range('D43').Value = "'39.28"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('D44').Value = "'0.265"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.52%  '
$ws.Range('D45').Value = '2.677.19'
$ws.Range('E45').Value = '  -0.92%  '
$ws.Range('D46').Value = "'132.23"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.09%  '
$ws.Range('D47').Value = "'0.0334"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.74%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value = "'1.00"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').Value = "'341.76"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.01%  '
$ws.Range('E50').Value = '  -1.76%  '
$ws.Range('D51').Value = "'21.40"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.37%  '
